# Applies the "Updated cryptos list on Mon Aug 28 14:29:09 UTC 2023 with GitHub
# Actions" commit: refreshed Price (D) / Volume(1h) (E) figures for most rows,
# plus three row swaps (Polkadot <-> WrappedEther, Frax <-> Cronos, and
# EnergySwap -> Mantle) in columns B/C/D/E on the single worksheet.
#
# Column D ("Price") cells are plain numeric-looking text (e.g. "219.20",
# "26.232.31") stored as strings in the source file. Writing such a string
# straight into `.Value` lets Excel auto-coerce it to a real number (dropping
# meaningful trailing zeros / reparsing multi-dot values), so each Price write
# is wrapped: force text format, assign, then restore the default "Normal"
# cell style so no stray numeric format lingers on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,4).NumberFormat = "@"
$ws.Cells.Item(2,4).Value = '26.232.31'
$ws.Cells.Item(2,4).Style = "Normal"
$ws.Cells.Item(2,5).Value = '  -0.64%  '

# Row 3
$ws.Cells.Item(3,4).NumberFormat = "@"
$ws.Cells.Item(3,4).Value = '1.656.20'
$ws.Cells.Item(3,4).Style = "Normal"
$ws.Cells.Item(3,5).Value = '  -0.90%  '

# Row 4
$ws.Cells.Item(4,5).Value = '  -0.68%  '

# Row 5
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = '219.20'
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(5,5).Value = '  -0.67%  '

# Row 6
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = '0.5234'
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Cells.Item(6,5).Value = '  -2.28%  '

# Row 7
$ws.Cells.Item(7,5).Value = '  -0.61%  '

# Row 8
$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = '0.2646'
$ws.Cells.Item(8,4).Style = "Normal"
$ws.Cells.Item(8,5).Value = '  -0.61%  '

# Row 9
$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value = '0.06324'
$ws.Cells.Item(9,4).Style = "Normal"
$ws.Cells.Item(9,5).Value = '  -1.16%  '

# Row 10
$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = '20.66'
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Cells.Item(10,5).Value = '  -1.70%  '

# Row 11
$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = '0.07774'
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Cells.Item(11,5).Value = '  -1.11%  '

# Row 12
$ws.Cells.Item(12,2).Value = 'WrappedEther'
$ws.Cells.Item(12,3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value = '1.738.83'
$ws.Cells.Item(12,4).Style = "Normal"
$ws.Cells.Item(12,5).Value = '  +3.95%  '

# Row 13
$ws.Cells.Item(13,2).Value = 'Polkadot'
$ws.Cells.Item(13,3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(13,4).NumberFormat = "@"
$ws.Cells.Item(13,4).Value = '4.559'
$ws.Cells.Item(13,4).Style = "Normal"
$ws.Cells.Item(13,5).Value = '  -0.09%  '

# Row 14
$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value = '1.884.04'
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Cells.Item(14,5).Value = '  -0.85%  '

# Row 15
$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value = '0.5640'
$ws.Cells.Item(15,4).Style = "Normal"
$ws.Cells.Item(15,5).Value = '  +1.40%  '

# Row 16
$ws.Cells.Item(16,4).NumberFormat = "@"
$ws.Cells.Item(16,4).Value = '0.0₅8093'
$ws.Cells.Item(16,4).Style = "Normal"
$ws.Cells.Item(16,5).Value = '  -1.12%  '

# Row 17
$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value = '65.36'
$ws.Cells.Item(17,4).Style = "Normal"
$ws.Cells.Item(17,5).Value = '  -1.34%  '

# Row 18
$ws.Cells.Item(18,4).NumberFormat = "@"
$ws.Cells.Item(18,4).Value = '26.224.18'
$ws.Cells.Item(18,4).Style = "Normal"
$ws.Cells.Item(18,5).Value = '  -0.75%  '

# Row 19
$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value = '1.004'
$ws.Cells.Item(19,4).Style = "Normal"
$ws.Cells.Item(19,5).Value = '  -0.59%  '

# Row 20
$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value = '4.727'
$ws.Cells.Item(20,4).Style = "Normal"
$ws.Cells.Item(20,5).Value = '  +0.78%  '

# Row 21
$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value = '193.00'
$ws.Cells.Item(21,4).Style = "Normal"
$ws.Cells.Item(21,5).Value = '  -1.74%  '

# Row 22
$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = '10.29'
$ws.Cells.Item(22,4).Style = "Normal"
$ws.Cells.Item(22,5).Value = '  -0.13%  '

# Row 23
$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value = '6.029'
$ws.Cells.Item(23,4).Style = "Normal"

# Row 24
$ws.Cells.Item(24,5).Value = '  -0.64%  '

# Row 25
$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value = '144.51'
$ws.Cells.Item(25,4).Style = "Normal"
$ws.Cells.Item(25,5).Value = '  -1.19%  '

# Row 26
$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value = '0.1204'
$ws.Cells.Item(26,4).Style = "Normal"
$ws.Cells.Item(26,5).Value = '  -2.07%  '

# Row 27
$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value = '7.270'
$ws.Cells.Item(27,4).Style = "Normal"
$ws.Cells.Item(27,5).Value = '  +0.23%  '

# Row 29
$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value = '1.503'
$ws.Cells.Item(29,4).Style = "Normal"
$ws.Cells.Item(29,5).Value = '  -0.27%  '

# Row 30
$ws.Cells.Item(30,4).NumberFormat = "@"
$ws.Cells.Item(30,4).Value = '0.05592'
$ws.Cells.Item(30,4).Style = "Normal"
$ws.Cells.Item(30,5).Value = '  -5.06%  '

# Row 31
$ws.Cells.Item(31,4).NumberFormat = "@"
$ws.Cells.Item(31,4).Value = '1.277'
$ws.Cells.Item(31,4).Style = "Normal"
$ws.Cells.Item(31,5).Value = '  -1.13%  '

# Row 32
$ws.Cells.Item(32,4).NumberFormat = "@"
$ws.Cells.Item(32,4).Value = '3.490'
$ws.Cells.Item(32,4).Style = "Normal"
$ws.Cells.Item(32,5).Value = '  -2.56%  '

# Row 33
$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(33,4).Value = '3.381'
$ws.Cells.Item(33,4).Style = "Normal"
$ws.Cells.Item(33,5).Value = '  +2.29%  '

# Row 34
$ws.Cells.Item(34,5).Value = '  -1.48%  '

# Row 35
$ws.Cells.Item(35,4).NumberFormat = "@"
$ws.Cells.Item(35,4).Value = '2.795'
$ws.Cells.Item(35,4).Style = "Normal"
$ws.Cells.Item(35,5).Value = '  -2.04%  '

# Row 36
$ws.Cells.Item(36,4).NumberFormat = "@"
$ws.Cells.Item(36,4).Value = '0.9464'
$ws.Cells.Item(36,4).Style = "Normal"
$ws.Cells.Item(36,5).Value = '  -2.57%  '

# Row 37
$ws.Cells.Item(37,4).NumberFormat = "@"
$ws.Cells.Item(37,4).Value = '2.401'
$ws.Cells.Item(37,4).Style = "Normal"
$ws.Cells.Item(37,5).Value = '  -1.24%  '

# Row 38
$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(38,4).Value = '0.5746'
$ws.Cells.Item(38,4).Style = "Normal"
$ws.Cells.Item(38,5).Value = '  -1.30%  '

# Row 39
$ws.Cells.Item(39,5).Value = '  -0.56%  '

# Row 40
$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = '5.953'
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Cells.Item(40,5).Value = '  +1.11%  '

# Row 41
$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value = '2.586'
$ws.Cells.Item(41,4).Style = "Normal"
$ws.Cells.Item(41,5).Value = '  -0.49%  '

# Row 42
$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = '0.8489'
$ws.Cells.Item(42,4).Style = "Normal"
$ws.Cells.Item(42,5).Value = '  -2.17%  '

# Row 43
$ws.Cells.Item(43,5).Value = '  -0.67%  '

# Row 44
$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value = '1.036.81'
$ws.Cells.Item(44,4).Style = "Normal"

# Row 45
$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = '102.48'
$ws.Cells.Item(45,4).Style = "Normal"
$ws.Cells.Item(45,5).Value = '  -1.70%  '

# Row 46
$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value = '1.794.68'
$ws.Cells.Item(46,4).Style = "Normal"
$ws.Cells.Item(46,5).Value = '  -0.83%  '

# Row 47
$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,4).Value = '58.32'
$ws.Cells.Item(47,4).Style = "Normal"
$ws.Cells.Item(47,5).Value = '  +0.24%  '

# Row 48
$ws.Cells.Item(48,5).Value = '  -1.18%  '

# Row 49
$ws.Cells.Item(49,2).Value = 'Cronos'
$ws.Cells.Item(49,3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,4).Value = '0.05319'
$ws.Cells.Item(49,4).Style = "Normal"
$ws.Cells.Item(49,5).Value = '  +2.91%  '

# Row 50
$ws.Cells.Item(50,2).Value = 'Frax'
$ws.Cells.Item(50,3).Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value = '0.9987'
$ws.Cells.Item(50,4).Style = "Normal"
$ws.Cells.Item(50,5).Value = '  -1.25%  '

# Row 51
$ws.Cells.Item(51,2).Value = 'Mantle'
$ws.Cells.Item(51,3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(51,4).NumberFormat = "@"
$ws.Cells.Item(51,4).Value = '0.4354'
$ws.Cells.Item(51,4).Style = "Normal"
$ws.Cells.Item(51,5).Value = '  -1.07%  '
